$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the execution date (column C) for all test rows from 17-12-2024 to 18-12-2024
$ws.Range("C2:C7").Value = "18-12-2024"

# Rename test case method names (column A)
$ws.Range("A4").Value = "verifyCustomerPackageUpgradeAndLicensePurchase"
$ws.Range("A6").Value = "verifyCustomerReceiptPageWithProratedAndRecurringOrderDetails"
$ws.Range("A7").Value = "verifyCustomerReceivedSubscriptionUpgradeReceipt"

# Re-fit column A width to the new (longer) content
$ws.Columns.Item(1).ColumnWidth = 60.65
